$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 461, shifting existing rows 461:559 down to 462:560
$ws.Rows("461:461").Insert()

# Populate the newly inserted row 461 with the new data record
$ws.Range("A461").Value = 10
$ws.Range("B461").Value = "Vega Modelo de Temuco"
$ws.Range("C461").Value = "La Araucanía"
$ws.Range("D461").Value = 44798
$ws.Range("D461").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E461").Value = 9
$ws.Range("F461").Value = 100112043
$ws.Range("G461").Value = "Pepino ensalada"
$ws.Range("H461").Value = "Sin especificar"
$ws.Range("I461").Value = "Primera"
$ws.Range("J461").Value = 700
$ws.Range("K461").Value = 24000
$ws.Range("L461").Value = 25000
$ws.Range("M461").Value = 24571
$ws.Range("N461").Value = "$/caja 60 unidades"
$ws.Range("O461").Value = "Región de Arica y Parinacota"
$ws.Range("P461").Value = 410
$ws.Range("Q461").Value = 60
$ws.Range("R461").Value = "Hortaliza"
